# --------------------------------------------------------------------
# Commit: "update import mockData and monthly used KwH stored per month"
#
# 1) Insert a new worksheet "MonthlyUsed" between "Client" and "Device"
#    holding per-device-per-month kWh usage (id, overallDeviceId, date,
#    monthlyKwHUsed).
# 2) Move the selection/cursor on the "Client" sheet to B6:B14 and on
#    the "Price" sheet to D2:D13 (cosmetic, reflects what the author was
#    last looking at).
# --------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- 1. Create & populate the new "MonthlyUsed" sheet -----------------
$deviceSheet = $wb.Worksheets.Item("Device")
$mu = $wb.Worksheets.Add($deviceSheet)
$mu.Name = "MonthlyUsed"

# Header row
$mu.Cells.Item(1, 1).Value = "id"
$mu.Cells.Item(1, 2).Value = "overallDeviceId"
$mu.Cells.Item(1, 3).Value = "date"
$mu.Cells.Item(1, 4).Value = "monthlyKwHUsed"

$monthlyData = @(
    @(1,1,45322,237),
    @(2,1,45351,216),
    @(3,1,45382,188),
    @(4,1,45412,253),
    @(5,1,45443,247),
    @(6,1,45473,177),
    @(7,1,45504,212),
    @(8,1,45535,229),
    @(9,1,45565,244),
    @(10,1,45596,220),
    @(11,1,45626,188),
    @(12,1,45657,184),
    @(13,2,45322,248),
    @(14,2,45351,213),
    @(15,2,45382,183),
    @(16,2,45412,246),
    @(17,2,45443,231),
    @(18,2,45473,203),
    @(19,2,45504,232),
    @(20,2,45535,222),
    @(21,2,45565,220),
    @(22,2,45596,200),
    @(23,2,45626,197),
    @(24,2,45657,215),
    @(25,3,45322,247),
    @(26,3,45351,182),
    @(27,3,45382,261),
    @(28,3,45412,215),
    @(29,3,45443,239),
    @(30,3,45473,226),
    @(31,3,45504,177),
    @(32,3,45535,206),
    @(33,3,45565,245),
    @(34,3,45596,187),
    @(35,3,45626,184),
    @(36,3,45657,204),
    @(37,4,45322,221),
    @(38,4,45351,238),
    @(39,4,45382,213),
    @(40,4,45412,245),
    @(41,4,45443,185),
    @(42,4,45473,239),
    @(43,4,45504,229),
    @(44,4,45535,253),
    @(45,4,45565,178),
    @(46,4,45596,256),
    @(47,4,45626,184),
    @(48,4,45657,207),
    @(49,5,45322,206),
    @(50,5,45351,217),
    @(51,5,45382,253),
    @(52,5,45412,211),
    @(53,5,45443,240),
    @(54,5,45473,223),
    @(55,5,45504,196),
    @(56,5,45535,258),
    @(57,5,45565,239),
    @(58,5,45596,206),
    @(59,5,45626,217),
    @(60,5,45657,224),
    @(61,6,45322,232),
    @(62,6,45351,225),
    @(63,6,45382,217),
    @(64,6,45412,216),
    @(65,6,45443,185),
    @(66,6,45473,205),
    @(67,6,45504,188),
    @(68,6,45535,261),
    @(69,6,45565,255),
    @(70,6,45596,254),
    @(71,6,45626,177),
    @(72,6,45657,208),
    @(73,7,45322,181),
    @(74,7,45351,244),
    @(75,7,45382,211),
    @(76,7,45412,210),
    @(77,7,45443,225),
    @(78,7,45473,202),
    @(79,7,45504,223),
    @(80,7,45535,221),
    @(81,7,45565,243),
    @(82,7,45596,185),
    @(83,7,45626,229),
    @(84,7,45657,231),
    @(85,8,45322,182),
    @(86,8,45351,194),
    @(87,8,45382,202),
    @(88,8,45412,198),
    @(89,8,45443,249),
    @(90,8,45473,232),
    @(91,8,45504,185),
    @(92,8,45535,230),
    @(93,8,45565,208),
    @(94,8,45596,199),
    @(95,8,45626,191),
    @(96,8,45657,237),
    @(97,9,45322,227),
    @(98,9,45351,178),
    @(99,9,45382,258),
    @(100,9,45412,248),
    @(101,9,45443,193),
    @(102,9,45473,229),
    @(103,9,45504,199),
    @(104,9,45535,222),
    @(105,9,45565,261),
    @(106,9,45596,211),
    @(107,9,45626,192),
    @(108,9,45657,256)
)

# Pre-format column C (date) to match the date style used by the other
# sheets (yyyy-mm-dd) before writing the serial date values in, so the
# existing number-format style gets reused.
$mu.Range("C2:C109").NumberFormat = "yyyy\-mm\-dd;@"

$r = 2
foreach ($row in $monthlyData) {
    $mu.Cells.Item($r, 1).Value = $row[0]
    $mu.Cells.Item($r, 2).Value = $row[1]
    $mu.Cells.Item($r, 3).Value = $row[2]
    $mu.Cells.Item($r, 4).Value = $row[3]
    $r++
}

$mu.Columns.Item(3).ColumnWidth = 10.33203125

# Cursor position on this new sheet
$mu.Range("F7").Select()

# ---- 2. Update selection on existing sheets ---------------------------
$clientSheet = $wb.Worksheets.Item("Client")
$clientSheet.Activate()
$clientSheet.Range("B6:B14").Select()

$priceSheet = $wb.Worksheets.Item("Price")
$priceSheet.Activate()
$priceSheet.Range("D2:D13").Select()

# ---- 3. Leave "MonthlyUsed" as the active tab --------------------------
$mu.Activate()
$mu.Range("F7").Select()
